$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.599.09"
$ws.Range("E2").Value = "  -0.77%  "
$ws.Range("D3").Value = "1.865.43"
$ws.Range("E3").Value = "  -1.08%  "
$ws.Range("D4").Value = "'1.016"
$ws.Range("E4").Value = "  -0.20%  "
$ws.Range("D5").Value = "'335.08"
$ws.Range("E5").Value = "  +0.32%  "
$ws.Range("E6").Value = "  -0.21%  "
$ws.Range("D7").Value = "'0.4658"
$ws.Range("E7").Value = "  -0.45%  "
$ws.Range("D8").Value = "'0.3922"
$ws.Range("E8").Value = "  +0.22%  "
$ws.Range("D9").Value = "'46.01"
$ws.Range("E9").Value = "  -3.02%  "
$ws.Range("D10").Value = "'0.07972"
$ws.Range("E10").Value = "  -1.06%  "
$ws.Range("D11").Value = "'1.000"
$ws.Range("E11").Value = "  -1.66%  "
$ws.Range("D12").Value = "'21.71"
$ws.Range("E12").Value = "  -1.06%  "
$ws.Range("D13").Value = "1.879.21"
$ws.Range("E13").Value = "  -1.63%  "
$ws.Range("D14").Value = "'5.958"
$ws.Range("E14").Value = "  +0.08%  "
$ws.Range("D15").Value = "'7.213"
$ws.Range("E15").Value = "  +1.94%  "
$ws.Range("E16").Value = "  -0.30%  "
$ws.Range("D17").Value = "'88.54"
$ws.Range("E17").Value = "  +1.49%  "
$ws.Range("D18").Value = "'0.06736"
$ws.Range("E18").Value = "  -0.40%  "
$ws.Range("D19").Value = "'0.00001043"
$ws.Range("E19").Value = "  -0.47%  "
$ws.Range("D20").Value = "'17.20"
$ws.Range("E20").Value = "  +0.10%  "
$ws.Range("D21").Value = "'1.014"
$ws.Range("E21").Value = "  -0.16%  "
$ws.Range("D22").Value = "27.616.36"
$ws.Range("E22").Value = "  -0.85%  "
$ws.Range("D23").Value = "'5.464"
$ws.Range("E23").Value = "  -0.59%  "
$ws.Range("D24").Value = "'10.92"
$ws.Range("E24").Value = "  -0.21%  "
$ws.Range("D25").Value = "'2.311"
$ws.Range("E25").Value = "  -1.36%  "
$ws.Range("D26").Value = "2.094.00"
$ws.Range("E26").Value = "  -1.83%  "
$ws.Range("D27").Value = "'159.25"
$ws.Range("E27").Value = "  -0.62%  "
$ws.Range("D28").Value = "'19.66"
$ws.Range("E28").Value = "  -2.02%  "
$ws.Range("D29").Value = "'2.134"
$ws.Range("E29").Value = "  +2.26%  "
$ws.Range("D30").Value = "'5.428"
$ws.Range("E30").Value = "  -1.30%  "
$ws.Range("D31").Value = "'121.74"
$ws.Range("E31").Value = "  -0.06%  "
$ws.Range("D32").Value = "'0.9752"
$ws.Range("E32").Value = "  +0.32%  "
$ws.Range("E33").Value = "  -0.17%  "
$ws.Range("D34").Value = "'3.629"
$ws.Range("E34").Value = "  -0.36%  "
$ws.Range("D35").Value = "'5.303"
$ws.Range("E35").Value = "  -0.90%  "
$ws.Range("D36").Value = "'1.337"
$ws.Range("E36").Value = "  -5.33%  "
$ws.Range("D37").Value = "'0.06026"
$ws.Range("E37").Value = "  -1.51%  "
$ws.Range("D38").Value = "'0.02232"
$ws.Range("E38").Value = "  -1.02%  "
$ws.Range("D39").Value = "'1.197"
$ws.Range("E39").Value = "  -1.33%  "
$ws.Range("D40").Value = "'8.288"
$ws.Range("E40").Value = "  +3.58%  "
$ws.Range("D41").Value = "'1.013"
$ws.Range("E41").Value = "  -0.29%  "
$ws.Range("D42").Value = "'0.5939"
$ws.Range("E42").Value = "  -0.65%  "
$ws.Range("D43").Value = "'0.1872"
$ws.Range("E43").Value = "  -0.84%  "
$ws.Range("D44").Value = "'10.30"
$ws.Range("E44").Value = "  +0.35%  "
$ws.Range("D45").Value = "'1.252"
$ws.Range("E45").Value = "  -0.96%  "
$ws.Range("D46").Value = "'0.5623"
$ws.Range("E46").Value = "  -1.00%  "
$ws.Range("D47").Value = "'12.19"
$ws.Range("E47").Value = "  -0.06%  "
$ws.Range("D48").Value = "'1.924"
$ws.Range("E48").Value = "  -0.24%  "
$ws.Range("D49").Value = "'0.06749"
$ws.Range("E49").Value = "  -2.48%  "
$ws.Range("D50").Value = "'111.72"
$ws.Range("E50").Value = "  -1.99%  "
$ws.Range("D51").Value = "'1.051"
$ws.Range("E51").Value = "  -1.71%  "
